$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 314, shifting existing rows 314:341 down to 315:342
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new record.
$ws.Cells.Item(314, 1).Value = 7
$ws.Cells.Item(314, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(314, 3).Value = "Ñuble"
$ws.Cells.Item(314, 4).Value = 45124
$ws.Cells.Item(314, 5).Value = 16
$ws.Cells.Item(314, 6).Value = 100112045
$ws.Cells.Item(314, 7).Value = "Zapallo"
$ws.Cells.Item(314, 8).Value = "Paine"
$ws.Cells.Item(314, 9).Value = "1a (guarda)"
$ws.Cells.Item(314, 10).Value = 200
$ws.Cells.Item(314, 11).Value = 350
$ws.Cells.Item(314, 12).Value = 350
$ws.Cells.Item(314, 13).Value = 350
$ws.Cells.Item(314, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(314, 15).Value = "Región del Maule"
$ws.Cells.Item(314, 16).Value = 350
$ws.Cells.Item(314, 17).Value = 1
$ws.Cells.Item(314, 18).Value = "Hortaliza"
